$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B labels (only rows where the label actually changes)
$ws.Range("B2").Value = "<arl>"
$ws.Range("B6").Value = "<seven>"
$ws.Range("B10").Value = "<lima>"
$ws.Range("B14").Value = "<are>"

# Update column C counts for rows 2-18
$ws.Range("C2").Value = 20
$ws.Range("C3").Value = 20
$ws.Range("C4").Value = 12
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 14
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 9
$ws.Range("C11").Value = 12
$ws.Range("C12").Value = 11
$ws.Range("C13").Value = 14
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 12
$ws.Range("C16").Value = 8
$ws.Range("C17").Value = 16
$ws.Range("C18").Value = 11
